$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1.2868788109663976,
    1.8627268598142788,
    1.8218534950195415,
    -4.8642896305829559,
    0.98234342384342854,
    -5.446363014985721,
    3.4843510502064023,
    -1.3857011763675899,
    3.5271551948617907,
    5.1530410122077743,
    -0.82341838989681881,
    -0.30761343189010987,
    -2.4200562058021693,
    -0.7311102773218765,
    -2.5445647622976608,
    -0.6428806835155384
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}

$ws.Range("A1:P1").Select()
